$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the May 1 2024 GitHub Actions run.
# D-column (Price) values are forced to text via NumberFormat so Excel
# does not silently reinterpret dotted/zero-padded strings as numbers;
# the format is reset back to Normal afterwards to avoid leaving a visible
# style change on the cell.

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '58.268.60'
$r.Style = 'Normal'
$ws.Range('E2').Value = '  -4.86%  '
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '2.931.21'
$r.Style = 'Normal'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('E4').Value = '  +0.10%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '556.71'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  -2.18%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '124.42'
$r.Style = 'Normal'
$ws.Range('E6').Value = '  -3.50%  '
$ws.Range('E7').Value = '  +0.15%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '2.921.32'
$r.Style = 'Normal'
$ws.Range('E8').Value = '  -2.79%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '0.496'
$r.Style = 'Normal'
$ws.Range('E9').Value = '  +0.16%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.127'
$r.Style = 'Normal'
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '4.78'
$r.Style = 'Normal'
$ws.Range('E11').Value = '  -7.23%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.440'
$r.Style = 'Normal'
$ws.Range('E12').Value = '  +2.36%  '
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '0.0000215'
$r.Style = 'Normal'
$ws.Range('E13').Value = '  -3.88%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '32.28'
$r.Style = 'Normal'
$ws.Range('E14').Value = '  -1.77%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '0.120'
$r.Style = 'Normal'
$ws.Range('E15').Value = '  +1.01%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '3.413.57'
$r.Style = 'Normal'
$ws.Range('E16').Value = '  -2.57%  '
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '2.927.81'
$r.Style = 'Normal'
$ws.Range('E17').Value = '  -2.56%  '
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '6.57'
$r.Style = 'Normal'
$ws.Range('E18').Value = '  +5.51%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '58.169.11'
$r.Style = 'Normal'
$ws.Range('E19').Value = '  -5.09%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '412.65'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  -6.12%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '12.94'
$r.Style = 'Normal'
$ws.Range('E21').Value = '  -1.99%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '0.667'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  +0.76%  '
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '6.90'
$r.Style = 'Normal'
$ws.Range('E23').Value = '  -3.41%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '12.95'
$r.Style = 'Normal'
$ws.Range('E24').Value = '  +3.01%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '77.47'
$r.Style = 'Normal'
$ws.Range('E25').Value = '  -1.86%  '
$ws.Range('E26').Value = '  +0.16%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '0.999'
$r.Style = 'Normal'
$ws.Range('E27').Value = '  +0.01%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '2.49'
$r.Style = 'Normal'
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '7.36'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  +0.70%  '
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '1.95'
$r.Style = 'Normal'
$ws.Range('E30').Value = '  +2.82%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '6.10'
$r.Style = 'Normal'
$ws.Range('E31').Value = '  -2.03%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '24.94'
$r.Style = 'Normal'
$ws.Range('E32').Value = '  -2.15%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '0.0992'
$r.Style = 'Normal'
$ws.Range('E33').Value = '  +4.98%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '0.925'
$r.Style = 'Normal'
$ws.Range('E34').Value = '  -3.14%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '2.04'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  -10.79%  '
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '5.42'
$r.Style = 'Normal'
$ws.Range('E36').Value = '  -2.63%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '48.17'
$r.Style = 'Normal'
$ws.Range('E37').Value = '  -3.80%  '
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '8.49'
$r.Style = 'Normal'
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.0₃0641'
$r.Style = 'Normal'
$ws.Range('E39').Value = '  -6.89%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '0.0350'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  -3.95%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.108'
$r.Style = 'Normal'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '366.31'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '2.635.49'
$r.Style = 'Normal'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '2.43'
$r.Style = 'Normal'
$ws.Range('E44').Value = '  -0.54%  '
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '0.233'
$r.Style = 'Normal'
$ws.Range('E46').Value = '  -1.19%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '118.98'
$r.Style = 'Normal'
$ws.Range('E47').Value = '  -1.07%  '
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '1.98'
$r.Style = 'Normal'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('E49').Value = '  +0.91%  '
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '23.03'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('E51').Value = '  -2.25%  '
